$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.002483333333333
$ws.Range("H2").Value = 9.007449999999999
$ws.Range("I2").Value = 0.1222246438870418
$ws.Range("J2").Value = 0.1222246438870418
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.918807333333334
$ws.Range("N2").Value = 29.756422
$ws.Range("O2").Value = 0.3718576623636295
$ws.Range("P2").Value = 0.3718576623636295
$ws.Range("Q2").Value = 29.78105370487777
$ws.Range("R2").Value = 268.0294833439
$ws.Range("S2").Value = 0.04545017035906244
$ws.Range("T2").Value = 0.04545017035906244
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.002483333333333
$ws.Range("H3").Value = 9.007449999999999
$ws.Range("I3").Value = 0.1222246438870418
$ws.Range("J3").Value = 0.1222246438870418
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.278342666666667
$ws.Range("N3").Value = 27.835028
$ws.Range("O3").Value = 0.3478465402831757
$ws.Range("P3").Value = 0.3478465402831757
$ws.Range("Q3").Value = 27.85806921762222
$ws.Range("R3").Value = 250.7226229586
$ws.Range("S3").Value = 0.04251541951345068
$ws.Range("T3").Value = 0.04251541951345068
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.002483333333333
$ws.Range("H4").Value = 9.007449999999999
$ws.Range("I4").Value = 0.1222246438870418
$ws.Range("J4").Value = 0.1222246438870418
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.476516666666666
$ws.Range("N4").Value = 22.42955
$ws.Range("O4").Value = 0.2802957973531948
$ws.Range("P4").Value = 0.2802957973531948
$ws.Range("Q4").Value = 22.44811668305555
$ws.Range("R4").Value = 202.0330501474999
$ws.Range("S4").Value = 0.03425905401452865
$ws.Range("T4").Value = 0.03425905401452865
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.971037
$ws.Range("H5").Value = 35.913111
$ws.Range("I5").Value = 0.4873151894099666
$ws.Range("J5").Value = 0.4873151894099665
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.918807333333334
$ws.Range("N5").Value = 29.756422
$ws.Range("O5").Value = 0.3718576623636295
$ws.Range("P5").Value = 0.3718576623636295
$ws.Range("Q5").Value = 118.7384095832047
$ws.Range("R5").Value = 1068.645686248842
$ws.Range("S5").Value = 0.1812118871682795
$ws.Range("T5").Value = 0.1812118871682795
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.971037
$ws.Range("H6").Value = 35.913111
$ws.Range("I6").Value = 0.4873151894099666
$ws.Range("J6").Value = 0.4873151894099665
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.278342666666667
$ws.Range("N6").Value = 27.835028
$ws.Range("O6").Value = 0.3478465402831757
$ws.Range("P6").Value = 0.3478465402831757
$ws.Range("Q6").Value = 111.0713833613453
$ws.Range("R6").Value = 999.642450252108
$ws.Range("S6").Value = 0.1695109026636973
$ws.Range("T6").Value = 0.1695109026636973
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.971037
$ws.Range("H7").Value = 35.913111
$ws.Range("I7").Value = 0.4873151894099666
$ws.Range("J7").Value = 0.4873151894099665
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.476516666666666
$ws.Range("N7").Value = 22.42955
$ws.Range("O7").Value = 0.2802957973531948
$ws.Range("P7").Value = 0.2802957973531948
$ws.Range("Q7").Value = 89.50165764778333
$ws.Range("R7").Value = 805.51491883005
$ws.Range("S7").Value = 0.1365923995779897
$ws.Range("T7").Value = 0.1365923995779897
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.591765666666666
$ws.Range("H8").Value = 28.775297
$ws.Range("I8").Value = 0.3904601667029916
$ws.Range("J8").Value = 0.3904601667029916
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.918807333333334
$ws.Range("N8").Value = 29.756422
$ws.Range("O8").Value = 0.3718576623636295
$ws.Range("P8").Value = 0.3718576623636295
$ws.Range("Q8").Value = 95.13887563414821
$ws.Range("R8").Value = 856.249880707334
$ws.Range("S8").Value = 0.1451956048362876
$ws.Range("T8").Value = 0.1451956048362876
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.591765666666666
$ws.Range("H9").Value = 28.775297
$ws.Range("I9").Value = 0.3904601667029916
$ws.Range("J9").Value = 0.3904601667029916
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.278342666666667
$ws.Range("N9").Value = 27.835028
$ws.Range("O9").Value = 0.3478465402831757
$ws.Range("P9").Value = 0.3478465402831757
$ws.Range("Q9").Value = 88.99568863370177
$ws.Range("R9").Value = 800.961197703316
$ws.Range("S9").Value = 0.1358202181060277
$ws.Range("T9").Value = 0.1358202181060277
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.591765666666666
$ws.Range("H10").Value = 28.775297
$ws.Range("I10").Value = 0.3904601667029916
$ws.Range("J10").Value = 0.3904601667029916
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.476516666666666
$ws.Range("N10").Value = 22.42955
$ws.Range("O10").Value = 0.2802957973531948
$ws.Range("P10").Value = 0.2802957973531948
$ws.Range("Q10").Value = 71.71299586959444
$ws.Range("R10").Value = 645.41696282635
$ws.Range("S10").Value = 0.1094443437606764
$ws.Range("T10").Value = 0.1094443437606764
